# Regenerate instance to have positive average demands during the last periods.
$wb = $excel.ActiveWorkbook

# Productdata!G2 : AverageDemand for product 1 -> 40 becomes 70
$wsProductdata = $wb.Worksheets.Item("Productdata")
$wsProductdata.Range("G2").Value = 70

# Workaround: the blank "t=s" (empty shared-string) cells in column H get
# corrupted into shared-string index 0 ("Name") by the COM runtime on
# save unless they are explicitly re-written as empty strings here.
$wsProductdata.Range("H2:H11").Value = ""

# ForecastedAverageDemand!B9:B11 (periods 7,8,9) -> 0 becomes 100
$wsAvgDemand = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvgDemand.Range("B9").Value = 100
$wsAvgDemand.Range("B10").Value = 100
$wsAvgDemand.Range("B11").Value = 100

# ForcastedStandardDeviation!B9:B11 (periods 7,8,9) -> 0 becomes positive decimals
$wsStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStdDev.Range("B9").Value = 10.23775
$wsStdDev.Range("B10").Value = 11.713975
$wsStdDev.Range("B11").Value = 13.0425775
